$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 9671.16
$ws.Range("J17").Value = 10011.625
$ws.Range("L17").Value = 30034.875
$ws.Range("N17").Value = -30370.875

$ws.Range("H137").Value = 2575.1948
$ws.Range("I137").Value = 1561.9
$ws.Range("J137").Value = 2930.7368
$ws.Range("K137").Value = 4685.700000000001
$ws.Range("L137").Value = 8792.2104
$ws.Range("M137").Value = -2135.700000000001
$ws.Range("N137").Value = -13892.2104

$ws.Range("H138").Value = 2976.4587
$ws.Range("I138").Value = 1789.9546
$ws.Range("J138").Value = 3390.7937
$ws.Range("K138").Value = 5369.8638
$ws.Range("L138").Value = 10172.3811
$ws.Range("M138").Value = -229.8638000000001
$ws.Range("N138").Value = -20452.3811

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10484.25
$ws.Range("I2").Value = 914.5
$ws.Range("J2").Value = 29623.75
$ws.Range("K2").Value = 914.5
$ws.Range("L2").Value = 29623.75
$ws.Range("M2").Value = -801.5
$ws.Range("N2").Value = -29849.75

$ws.Range("H32").Value = 4026.5881
$ws.Range("I32").Value = 3223.7344
$ws.Range("J32").Value = 16872.25
$ws.Range("K32").Value = 3223.7344
$ws.Range("L32").Value = 16872.25
$ws.Range("M32").Value = -2936.7344
$ws.Range("N32").Value = -17446.25

$ws.Range("H45").Value = 62502170
$ws.Range("I45").Value = 62502170
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 62502170
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -62501793

$ws.Range("H97").Value = 2035
$ws.Range("I97").Value = 2035
$ws.Range("K97").Value = 2035
$ws.Range("M97").Value = -1539

$ws.Range("H116").Value = 10484.25
$ws.Range("I116").Value = 914.5
$ws.Range("J116").Value = 29623.75
$ws.Range("K116").Value = 914.5
$ws.Range("L116").Value = 29623.75
$ws.Range("M116").Value = 1379.5
$ws.Range("N116").Value = -34211.75

$ws.Range("H132").Value = 2798.3655
$ws.Range("I132").Value = 2103.5952
$ws.Range("J132").Value = 5716.4
$ws.Range("K132").Value = 6310.785600000001
$ws.Range("L132").Value = 17149.2
$ws.Range("M132").Value = -3780.785600000001
$ws.Range("N132").Value = -22209.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10484.25
$ws.Range("I3").Value = 914.5
$ws.Range("J3").Value = 29623.75
$ws.Range("K3").Value = 914.5
$ws.Range("L3").Value = 29623.75
$ws.Range("M3").Value = -800.5
$ws.Range("N3").Value = -29851.75

$ws.Range("H13").Value = 70000
$ws.Range("J13").Value = 70000
$ws.Range("L13").Value = 70000
$ws.Range("N13").Value = -70336

$ws.Range("H105").Value = 9674.757
$ws.Range("J105").Value = 9841.429
$ws.Range("L105").Value = 9841.429
$ws.Range("N105").Value = -13335.429

$ws.Range("H106").Value = 33931.168
$ws.Range("J106").Value = 33931.168
$ws.Range("L106").Value = 33931.168
$ws.Range("N106").Value = -36455.168

$ws.Range("H107").Value = 1674.3334
$ws.Range("I107").Value = 1672.375
$ws.Range("K107").Value = 1672.375
$ws.Range("M107").Value = 247.625

$ws.Range("H132").Value = 69849.5
$ws.Range("J132").Value = 69849.5
$ws.Range("L132").Value = 69849.5
$ws.Range("N132").Value = -79969.5

$ws.Range("H134").Value = 1897.9231
$ws.Range("I134").Value = 1166.6086
$ws.Range("K134").Value = 3499.8258
$ws.Range("M134").Value = -964.8258000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 39093.45
$ws.Range("I31").Value = 2032.7142
$ws.Range("J31").Value = 136377.88
$ws.Range("K31").Value = 2032.7142
$ws.Range("L31").Value = 136377.88
$ws.Range("M31").Value = -1737.7142
$ws.Range("N31").Value = -136967.88

$ws.Range("H34").Value = 39093.45
$ws.Range("I34").Value = 2032.7142
$ws.Range("J34").Value = 136377.88
$ws.Range("K34").Value = 2032.7142
$ws.Range("L34").Value = 136377.88
$ws.Range("M34").Value = -1830.7142
$ws.Range("N34").Value = -136781.88

$ws.Range("H50").Value = 38000
$ws.Range("J50").Value = 38000
$ws.Range("L50").Value = 38000
$ws.Range("N50").Value = -39250

$ws.Range("H105").Value = 4328.0835
$ws.Range("I105").Value = 1587
$ws.Range("K105").Value = 1587
$ws.Range("M105").Value = 160

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 11001
$ws.Range("I80").Value = 6999
$ws.Range("J80").Value = 15003
$ws.Range("K80").Value = 20997
$ws.Range("L80").Value = 45009
$ws.Range("M80").Value = -20061
$ws.Range("N80").Value = -46881

$ws.Range("H83").Value = 11001
$ws.Range("I83").Value = 6999
$ws.Range("J83").Value = 15003
$ws.Range("K83").Value = 62991
$ws.Range("L83").Value = 135027
$ws.Range("M83").Value = -58311
$ws.Range("N83").Value = -144387

$ws.Range("H132").Value = 3286.1155
$ws.Range("I132").Value = 2963.923
$ws.Range("K132").Value = 26675.307
$ws.Range("M132").Value = -24145.307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6452.6924
$ws.Range("I122").Value = 5192.72
$ws.Range("K122").Value = 15578.16
$ws.Range("M122").Value = -13128.16

$ws.Range("H123").Value = 36530.4
$ws.Range("J123").Value = 36530.4
$ws.Range("L123").Value = 36530.4
$ws.Range("N123").Value = -41430.4

$ws.Range("H132").Value = 2591.96
$ws.Range("I132").Value = 2831.2307
$ws.Range("J132").Value = 2332.75
$ws.Range("K132").Value = 8493.6921
$ws.Range("L132").Value = 6998.25
$ws.Range("M132").Value = -5963.6921
$ws.Range("N132").Value = -12058.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6450.5806
$ws.Range("I46").Value = 5179.5557
$ws.Range("K46").Value = 5179.5557
$ws.Range("M46").Value = -4991.5557

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 12505.2
$ws.Range("J14").Value = 8763
$ws.Range("L14").Value = 8763
$ws.Range("N14").Value = -9099

$ws.Range("H132").Value = 4700.4287
$ws.Range("I132").Value = 2339.8
$ws.Range("J132").Value = 10602
$ws.Range("K132").Value = 7019.400000000001
$ws.Range("L132").Value = 31806
$ws.Range("M132").Value = -4489.400000000001
$ws.Range("N132").Value = -36866

$ws.Range("H136").Value = 2897.98
$ws.Range("I136").Value = 2023.6
$ws.Range("J136").Value = 4938.2
$ws.Range("K136").Value = 6070.799999999999
$ws.Range("L136").Value = 14814.6
$ws.Range("M136").Value = -3520.799999999999
$ws.Range("N136").Value = -19914.6
